$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before the existing row 411, shifting rows 411-420
# down to 415-424 (mirrors a weekly data refresh: a new reporting date's
# rows were inserted ahead of the previously-last date's rows).
$ws.Range("A411:T414").Insert()

# --- New row 411 (new date 2021-09-09 = serial 44448) ---
$ws.Range("A411").Value = 6
$ws.Range("B411").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C411").Value = "Metropolitana"
$ws.Range("D411").Value = 44448
$ws.Range("E411").Value = 13
$ws.Range("F411").Value = "Fruta"
$ws.Range("G411").Value = 100101
$ws.Range("H411").Value = "Berries"
$ws.Range("I411").Value = 100101007
$ws.Range("J411").Value = "Kiwi"
$ws.Range("K411").Value = "Hayward"
$ws.Range("L411").Value = "Especial"
$ws.Range("M411").Value = 20
$ws.Range("N411").Value = 320000
$ws.Range("O411").Value = 330000
$ws.Range("P411").Value = 325000
$ws.Range("Q411").Value = "$/bins (450 kilos)"
$ws.Range("R411").Value = "Región de O'Higgins"
$ws.Range("S411").Value = 722
$ws.Range("T411").Value = 450

# --- New row 412 ---
$ws.Range("A412").Value = 6
$ws.Range("B412").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C412").Value = "Metropolitana"
$ws.Range("D412").Value = 44448
$ws.Range("E412").Value = 13
$ws.Range("F412").Value = "Fruta"
$ws.Range("G412").Value = 100101
$ws.Range("H412").Value = "Berries"
$ws.Range("I412").Value = 100101007
$ws.Range("J412").Value = "Kiwi"
$ws.Range("K412").Value = "Hayward"
$ws.Range("L412").Value = "Extra (doble especial)"
$ws.Range("M412").Value = 15
$ws.Range("N412").Value = 420000
$ws.Range("O412").Value = 420000
$ws.Range("P412").Value = 420000
$ws.Range("Q412").Value = "$/bins (450 kilos)"
$ws.Range("R412").Value = "Región de O'Higgins"
$ws.Range("S412").Value = 933
$ws.Range("T412").Value = 450

# --- New row 413 ---
$ws.Range("A413").Value = 6
$ws.Range("B413").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C413").Value = "Metropolitana"
$ws.Range("D413").Value = 44448
$ws.Range("E413").Value = 13
$ws.Range("F413").Value = "Fruta"
$ws.Range("G413").Value = 100101
$ws.Range("H413").Value = "Berries"
$ws.Range("I413").Value = 100101007
$ws.Range("J413").Value = "Kiwi"
$ws.Range("K413").Value = "Hayward"
$ws.Range("L413").Value = "Primera"
$ws.Range("M413").Value = 30
$ws.Range("N413").Value = 250000
$ws.Range("O413").Value = 260000
$ws.Range("P413").Value = 255000
$ws.Range("Q413").Value = "$/bins (450 kilos)"
$ws.Range("R413").Value = "Región de O'Higgins"
$ws.Range("S413").Value = 567
$ws.Range("T413").Value = 450

# --- New row 414 ---
$ws.Range("A414").Value = 6
$ws.Range("B414").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C414").Value = "Metropolitana"
$ws.Range("D414").Value = 44448
$ws.Range("E414").Value = 13
$ws.Range("F414").Value = "Fruta"
$ws.Range("G414").Value = 100101
$ws.Range("H414").Value = "Berries"
$ws.Range("I414").Value = 100101007
$ws.Range("J414").Value = "Kiwi"
$ws.Range("K414").Value = "Hayward"
$ws.Range("L414").Value = "Segunda"
$ws.Range("M414").Value = 20
$ws.Range("N414").Value = 220000
$ws.Range("O414").Value = 220000
$ws.Range("P414").Value = 220000
$ws.Range("Q414").Value = "$/bins (450 kilos)"
$ws.Range("R414").Value = "Región de O'Higgins"
$ws.Range("S414").Value = 489
$ws.Range("T414").Value = 450
